# The deck currently has its (only reachable/applied) DrawingML colour
# theme set to the "Integral" palette. The authored change swaps the
# theme applied to the deck back to the stock "Office Theme" palette
# (the 12 dk/lt/accent/hyperlink colours of the default Office theme).
#
# Helper: convert an RRGGBB hex string into the OLE_COLOR / VBA RGB()
# integer layout (R + G*256 + B*65536) that ThemeColor.RGB expects.
function HexToRgbInt([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme colour scheme, in MsoThemeColorSchemeIndex order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink.
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = HexToRgbInt $officeThemeColors[$i - 1]
}
